# v2.0 update: add two new match rows for Varun Chakravarthy and keep the
# existing match row, matching the order: Oct 12 2020 (new), Oct 7 2020
# (existing, pushed down), Oct 26 2020 (new, appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data's playerName value ends with U+00A0 (NBSP), not a plain
# space - match it exactly so the new rows' F column is byte-identical to
# the existing F column value.
$playerName = "Varun Chakravarthy" + [char]0x00A0

# Insert a brand-new row above the existing data row (old row 2), so the
# existing "Oct 7 2020" record slides down to row 3.
$ws.Rows(2).Insert()

# Columns G:K hold numeric-looking values that must stay stored as TEXT
# (as in the original file), so force text format before writing them.
$ws.Range("G2:K4").NumberFormat = "@"

# New row 2: Oct 12 2020 vs Royal Challengers Bangalore
$ws.Range("A2").Value = " Oct 12 2020"
$ws.Range("B2").Value = " Sharjah"
$ws.Range("C2").Value = "RCB won by 82 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Royal Challengers Bangalore"
$ws.Range("F2").Value = $playerName
$ws.Range("G2").Value = "7"
$ws.Range("H2").Value = "10"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "70.00"

# New row 4: Oct 26 2020 vs Kings XI Punjab (appended after the existing
# "Oct 7 2020" row, which is now row 3)
$ws.Range("A4").Value = " Oct 26 2020"
$ws.Range("B4").Value = " Sharjah"
$ws.Range("C4").Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Kings XI Punjab"
$ws.Range("F4").Value = $playerName
$ws.Range("G4").Value = "2"
$ws.Range("H4").Value = "4"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "50.00"
